$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 211.4614666666667
$ws.Range("H2").Value = 634.3844
$ws.Range("I2").Value = 0.2421062275331183
$ws.Range("J2").Value = 0.2421062275331183
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9442423333333334
$ws.Range("N2").Value = 2.832727
$ws.Range("O2").Value = 0.006848500623481535
$ws.Range("P2").Value = 0.006848500623481536
$ws.Range("Q2").Value = 199.6708686954222
$ws.Range("R2").Value = 1797.0378182588
$ws.Range("S2").Value = 0.001658064650209323
$ws.Range("T2").Value = 0.001658064650209323

$ws.Range("G3").Value = 211.4614666666667
$ws.Range("H3").Value = 634.3844
$ws.Range("I3").Value = 0.2421062275331183
$ws.Range("J3").Value = 0.2421062275331183
$ws.Range("M3").Value = 82.477727
$ws.Range("N3").Value = 247.433181
$ws.Range("O3").Value = 0.5982031781913751
$ws.Range("P3").Value = 0.5982031781913751
$ws.Range("Q3").Value = 17440.86111875293
$ws.Range("R3").Value = 156967.7500687764
$ws.Range("S3").Value = 0.1448287147702356
$ws.Range("T3").Value = 0.1448287147702356

$ws.Range("G4").Value = 211.4614666666667
$ws.Range("H4").Value = 634.3844
$ws.Range("I4").Value = 0.2421062275331183
$ws.Range("J4").Value = 0.2421062275331183
$ws.Range("M4").Value = 9.766934000000001
$ws.Range("N4").Value = 29.300802
$ws.Range("O4").Value = 0.07083865150630789
$ws.Range("P4").Value = 0.07083865150630789
$ws.Range("Q4").Value = 2065.330188476533
$ws.Range("R4").Value = 18587.9716962888
$ws.Range("S4").Value = 0.01715047867972545
$ws.Range("T4").Value = 0.01715047867972545

$ws.Range("G5").Value = 211.4614666666667
$ws.Range("H5").Value = 634.3844
$ws.Range("I5").Value = 0.2421062275331183
$ws.Range("J5").Value = 0.2421062275331183
$ws.Range("M5").Value = 44.68687199999999
$ws.Range("N5").Value = 134.060616
$ws.Range("O5").Value = 0.3241096696788354
$ws.Range("P5").Value = 0.3241096696788355
$ws.Range("Q5").Value = 9449.551493865598
$ws.Range("R5").Value = 85045.9634447904
$ws.Range("S5").Value = 0.07846896943294794
$ws.Range("T5").Value = 0.07846896943294795

$ws.Range("I6").Value = 0.08842543241393927
$ws.Range("J6").Value = 0.08842543241393927
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9442423333333334
$ws.Range("N6").Value = 2.832727
$ws.Range("O6").Value = 0.006848500623481535
$ws.Range("P6").Value = 0.006848500623481536
$ws.Range("Q6").Value = 72.92659542367367
$ws.Range("R6").Value = 656.3393588130631
$ws.Range("S6").Value = 0.0006055816290184875
$ws.Range("T6").Value = 0.0006055816290184875

$ws.Range("I7").Value = 0.08842543241393927
$ws.Range("J7").Value = 0.08842543241393927
$ws.Range("M7").Value = 82.477727
$ws.Range("N7").Value = 247.433181
$ws.Range("O7").Value = 0.5982031781913751
$ws.Range("P7").Value = 0.5982031781913751
$ws.Range("Q7").Value = 6369.995938606021
$ws.Range("R7").Value = 57329.96344745419
$ws.Range("S7").Value = 0.05289637470296511
$ws.Range("T7").Value = 0.05289637470296511

$ws.Range("I8").Value = 0.08842543241393927
$ws.Range("J8").Value = 0.08842543241393927
$ws.Range("M8").Value = 9.766934000000001
$ws.Range("N8").Value = 29.300802
$ws.Range("O8").Value = 0.07083865150630789
$ws.Range("P8").Value = 0.07083865150630789
$ws.Range("Q8").Value = 754.3288615680821
$ws.Range("R8").Value = 6788.959754112738
$ws.Range("S8").Value = 0.006263938391065625
$ws.Range("T8").Value = 0.006263938391065625

$ws.Range("I9").Value = 0.08842543241393927
$ws.Range("J9").Value = 0.08842543241393927
$ws.Range("M9").Value = 44.68687199999999
$ws.Range("N9").Value = 134.060616
$ws.Range("O9").Value = 0.3241096696788354
$ws.Range("P9").Value = 0.3241096696788355
$ws.Range("Q9").Value = 3451.297744286856
$ws.Range("R9").Value = 31061.6796985817
$ws.Range("S9").Value = 0.02865953769089004
$ws.Range("T9").Value = 0.02865953769089005

$ws.Range("G10").Value = 174.3107043333333
$ws.Range("H10").Value = 522.932113
$ws.Range("I10").Value = 0.199571617988009
$ws.Range("J10").Value = 0.199571617988009
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.9442423333333334
$ws.Range("N10").Value = 2.832727
$ws.Range("O10").Value = 0.006848500623481535
$ws.Range("P10").Value = 0.006848500623481536
$ws.Range("Q10").Value = 164.5915461846834
$ws.Range("R10").Value = 1481.323915662151
$ws.Range("S10").Value = 0.001366766350220099
$ws.Range("T10").Value = 0.001366766350220099

$ws.Range("G11").Value = 174.3107043333333
$ws.Range("H11").Value = 522.932113
$ws.Range("I11").Value = 0.199571617988009
$ws.Range("J11").Value = 0.199571617988009
$ws.Range("M11").Value = 82.477727
$ws.Range("N11").Value = 247.433181
$ws.Range("O11").Value = 0.5982031781913751
$ws.Range("P11").Value = 0.5982031781913751
$ws.Range("Q11").Value = 14376.75068518238
$ws.Range("R11").Value = 129390.7561666414
$ws.Range("S11").Value = 0.119384376157222
$ws.Range("T11").Value = 0.119384376157222

$ws.Range("G12").Value = 174.3107043333333
$ws.Range("H12").Value = 522.932113
$ws.Range("I12").Value = 0.199571617988009
$ws.Range("J12").Value = 0.199571617988009
$ws.Range("M12").Value = 9.766934000000001
$ws.Range("N12").Value = 29.300802
$ws.Range("O12").Value = 0.07083865150630789
$ws.Range("P12").Value = 0.07083865150630789
$ws.Range("Q12").Value = 1702.481144717181
$ws.Range("R12").Value = 15322.33030245463
$ws.Range("S12").Value = 0.01413738429720258
$ws.Range("T12").Value = 0.01413738429720258

$ws.Range("G13").Value = 174.3107043333333
$ws.Range("H13").Value = 522.932113
$ws.Range("I13").Value = 0.199571617988009
$ws.Range("J13").Value = 0.199571617988009
$ws.Range("M13").Value = 44.68687199999999
$ws.Range("N13").Value = 134.060616
$ws.Range("O13").Value = 0.3241096696788354
$ws.Range("P13").Value = 0.3241096696788355
$ws.Range("Q13").Value = 7789.40013277351
$ws.Range("R13").Value = 70104.6011949616
$ws.Range("S13").Value = 0.06468309118336434
$ws.Range("T13").Value = 0.06468309118336434

$ws.Range("G14").Value = 28.53474833333333
$ws.Range("H14").Value = 85.60424499999999
$ws.Range("I14").Value = 0.03266997236655063
$ws.Range("J14").Value = 0.03266997236655063
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.9442423333333334
$ws.Range("N14").Value = 2.832727
$ws.Range("O14").Value = 0.006848500623481535
$ws.Range("P14").Value = 0.006848500623481536
$ws.Range("Q14").Value = 26.94371734734611
$ws.Range("R14").Value = 242.493456126115
$ws.Range("S14").Value = 0.0002237403261214465
$ws.Range("T14").Value = 0.0002237403261214465

$ws.Range("G15").Value = 28.53474833333333
$ws.Range("H15").Value = 85.60424499999999
$ws.Range("I15").Value = 0.03266997236655063
$ws.Range("J15").Value = 0.03266997236655063
$ws.Range("M15").Value = 82.477727
$ws.Range("N15").Value = 247.433181
$ws.Range("O15").Value = 0.5982031781913751
$ws.Range("P15").Value = 0.5982031781913751
$ws.Range("Q15").Value = 2353.481183050371
$ws.Range("R15").Value = 21181.33064745334
$ws.Range("S15").Value = 0.01954328130109499
$ws.Range("T15").Value = 0.01954328130109499

$ws.Range("G16").Value = 28.53474833333333
$ws.Range("H16").Value = 85.60424499999999
$ws.Range("I16").Value = 0.03266997236655063
$ws.Range("J16").Value = 0.03266997236655063
$ws.Range("M16").Value = 9.766934000000001
$ws.Range("N16").Value = 29.300802
$ws.Range("O16").Value = 0.07083865150630789
$ws.Range("P16").Value = 0.07083865150630789
$ws.Range("Q16").Value = 278.6970036782766
$ws.Range("R16").Value = 2508.27303310449
$ws.Range("S16").Value = 0.002314296787194789
$ws.Range("T16").Value = 0.002314296787194789

$ws.Range("G17").Value = 28.53474833333333
$ws.Range("H17").Value = 85.60424499999999
$ws.Range("I17").Value = 0.03266997236655063
$ws.Range("J17").Value = 0.03266997236655063
$ws.Range("M17").Value = 44.68687199999999
$ws.Range("N17").Value = 134.060616
$ws.Range("O17").Value = 0.3241096696788354
$ws.Range("P17").Value = 0.3241096696788355
$ws.Range("Q17").Value = 1275.12864632388
$ws.Range("R17").Value = 11476.15781691492
$ws.Range("S17").Value = 0.01058865395213941
$ws.Range("T17").Value = 0.01058865395213941

$ws.Range("G18").Value = 230.32901
$ws.Range("H18").Value = 690.98703
$ws.Range("I18").Value = 0.263708034289011
$ws.Range("J18").Value = 0.263708034289011
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.9442423333333334
$ws.Range("N18").Value = 2.832727
$ws.Range("O18").Value = 0.006848500623481535
$ws.Range("P18").Value = 0.006848500623481536
$ws.Range("Q18").Value = 217.4864018367567
$ws.Range("R18").Value = 1957.37761653081
$ws.Range("S18").Value = 0.001806004637245382
$ws.Range("T18").Value = 0.001806004637245382

$ws.Range("G19").Value = 230.32901
$ws.Range("H19").Value = 690.98703
$ws.Range("I19").Value = 0.263708034289011
$ws.Range("J19").Value = 0.263708034289011
$ws.Range("M19").Value = 82.477727
$ws.Range("N19").Value = 247.433181
$ws.Range("O19").Value = 0.5982031781913751
$ws.Range("P19").Value = 0.5982031781913751
$ws.Range("Q19").Value = 18997.01320696027
$ws.Range("R19").Value = 170973.1188626424
$ws.Range("S19").Value = 0.1577509842262865
$ws.Range("T19").Value = 0.1577509842262865

$ws.Range("G20").Value = 230.32901
$ws.Range("H20").Value = 690.98703
$ws.Range("I20").Value = 0.263708034289011
$ws.Range("J20").Value = 0.263708034289011
$ws.Range("M20").Value = 9.766934000000001
$ws.Range("N20").Value = 29.300802
$ws.Range("O20").Value = 0.07083865150630789
$ws.Range("P20").Value = 0.07083865150630789
$ws.Range("Q20").Value = 2249.60823895534
$ws.Range("R20").Value = 20246.47415059806
$ws.Range("S20").Value = 0.01868072154041274
$ws.Range("T20").Value = 0.01868072154041274

$ws.Range("G21").Value = 230.32901
$ws.Range("H21").Value = 690.98703
$ws.Range("I21").Value = 0.263708034289011
$ws.Range("J21").Value = 0.263708034289011
$ws.Range("M21").Value = 44.68687199999999
$ws.Range("N21").Value = 134.060616
$ws.Range("O21").Value = 0.3241096696788354
$ws.Range("P21").Value = 0.3241096696788355
$ws.Range("Q21").Value = 10292.68298775672
$ws.Range("R21").Value = 92634.14688981046
$ws.Range("S21").Value = 0.08547032388506635
$ws.Range("T21").Value = 0.08547032388506637

$ws.Range("G22").Value = 151.5554656666667
$ws.Range("H22").Value = 454.666397
$ws.Range("I22").Value = 0.1735187154093718
$ws.Range("J22").Value = 0.1735187154093718
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 0.9442423333333334
$ws.Range("N22").Value = 2.832727
$ws.Range("O22").Value = 0.006848500623481535
$ws.Range("P22").Value = 0.006848500623481536
$ws.Range("Q22").Value = 143.1050865305132
$ws.Range("R22").Value = 1287.945778774619
$ws.Range("S22").Value = 0.001188343030666798
$ws.Range("T22").Value = 0.001188343030666798

$ws.Range("G23").Value = 151.5554656666667
$ws.Range("H23").Value = 454.666397
$ws.Range("I23").Value = 0.1735187154093718
$ws.Range("J23").Value = 0.1735187154093718
$ws.Range("M23").Value = 82.477727
$ws.Range("N23").Value = 247.433181
$ws.Range("O23").Value = 0.5982031781913751
$ws.Range("P23").Value = 0.5982031781913751
$ws.Range("Q23").Value = 12499.95032261321
$ws.Range("R23").Value = 112499.5529035188
$ws.Range("S23").Value = 0.103799447033571
$ws.Range("T23").Value = 0.103799447033571

$ws.Range("G24").Value = 151.5554656666667
$ws.Range("H24").Value = 454.666397
$ws.Range("I24").Value = 0.1735187154093718
$ws.Range("J24").Value = 0.1735187154093718
$ws.Range("M24").Value = 9.766934000000001
$ws.Range("N24").Value = 29.300802
$ws.Range("O24").Value = 0.07083865150630789
$ws.Range("P24").Value = 0.07083865150630789
$ws.Range("Q24").Value = 1480.2322305056
$ws.Range("R24").Value = 13322.09007455039
$ws.Range("S24").Value = 0.01229183181070671
$ws.Range("T24").Value = 0.01229183181070671

$ws.Range("G25").Value = 151.5554656666667
$ws.Range("H25").Value = 454.666397
$ws.Range("I25").Value = 0.1735187154093718
$ws.Range("J25").Value = 0.1735187154093718
$ws.Range("M25").Value = 44.68687199999999
$ws.Range("N25").Value = 134.060616
$ws.Range("O25").Value = 0.3241096696788354
$ws.Range("P25").Value = 0.3241096696788355
$ws.Range("Q25").Value = 6772.539695146727
$ws.Range("R25").Value = 60952.85725632054
$ws.Range("S25").Value = 0.05623909353442735
$ws.Range("T25").Value = 0.05623909353442735
